# "Generate Report for Handback"
#
# The localization CI job re-ran after the de-de handback file caught up
# with en-US, so for both language sheets (and the roll-up Overview
# sheet) the Status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the "Latest Handback DateTime"
# stamps advance, and the stale "Error Detail" (version-mismatch
# warning) clears out now that everything is in sync. The Status /
# Error Detail columns are also widened so the new, longer text fits.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Column width as stored in the xlsx is ~ (Excel ColumnWidth + 5/6) rounded
# to the nearest 1/6 of a character; pick the ColumnWidth whose stored
# result lands closest to the desired on-disk width.
$wideStatusWidth = 29.166666666666668   # -> stored width ~29.98 (was ~17.22)
$narrowErrorWidth = 12.833333333333332  # -> stored width ~13.75 (was 40)

# ---------------------------------------------------------------------
# Overview sheet: roll-up Status columns for zh-cn (E) and de-de (F)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value2 = $newStatus
$wsOverview.Range("F2").Value2 = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $wideStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatusWidth

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value2 = $newStatus
$wsZhCn.Range("K2").Value2 = "2016-09-05 06:56:31"
$wsZhCn.Range("P2").Value2 = ""
$wsZhCn.Columns.Item(3).ColumnWidth = $wideStatusWidth
$wsZhCn.Columns.Item(16).ColumnWidth = $narrowErrorWidth

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value2 = $newStatus
$wsDeDe.Range("K2").Value2 = "2016-09-05 06:56:39"
$wsDeDe.Range("P2").Value2 = ""
$wsDeDe.Columns.Item(3).ColumnWidth = $wideStatusWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $narrowErrorWidth
